# Updates the cryptos list with new prices / volume figures, and reorders
# the Hedera / VeChain rows (37-38) to reflect the latest ranking.
# Price (column D) cells are plain text in this sheet (e.g. "27.369.34"),
# so force a text NumberFormat before writing to stop Excel from
# re-interpreting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37 / 38 swap: VeChain now ranks above Hedera ---
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02421"
$ws.Range("E37").Value = "  +0.73%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06525"
$ws.Range("E38").Value = "  +0.01%  "

# --- Price / Volume(1h) refresh for the remaining rows ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.369.34"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.790.64"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "340.26"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3938"
$ws.Range("E7").Value = "  +2.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3463"
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.18"
$ws.Range("E9").Value = "  -3.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.198"
$ws.Range("E10").Value = "  -3.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07496"
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9999"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.81"
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.509"
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.787.60"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.142"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001096"
$ws.Range("E17").Value = "  -2.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06688"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.75"
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9988"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.73"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.532"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.369.55"
$ws.Range("E23").Value = "  -0.93%  "
$ws.Range("E24").Value = "  -5.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.408"
$ws.Range("E25").Value = "  -2.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.22"
$ws.Range("E26").Value = "  -4.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.505"
$ws.Range("E27").Value = "  -6.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.460"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.12"
$ws.Range("E29").Value = "  +3.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.990.22"
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "136.12"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.028"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.016"
$ws.Range("E33").Value = "  -4.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08843"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.02"
$ws.Range("E35").Value = "  -6.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.625"
$ws.Range("E36").Value = "  -4.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.422"
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6825"
$ws.Range("E40").Value = "  -2.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2217"
$ws.Range("E41").Value = "  -1.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.253"
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.382"
$ws.Range("E43").Value = "  -8.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.44"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9986"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6392"
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.874"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.136"
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.35"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07163"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.177"
$ws.Range("E51").Value = "  +3.01%  "
